# Update the crypto price/volume table on Sheet1.
# Each cell is set with a leading apostrophe so Excel stores the value as
# literal text (matching the workbook's original inline-string cells)
# instead of auto-converting numeric-looking strings (e.g. "212.78") into
# actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.679.37"
$ws.Range("E2").Value = "'  +1.81%  "
$ws.Range("D3").Value = "'1.634.93"
$ws.Range("E3").Value = "'  +1.95%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'212.78"
$ws.Range("D6").Value = "'0.495"
$ws.Range("E6").Value = "'  +2.19%  "
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E8").Value = "'  +1.38%  "
$ws.Range("E9").Value = "'  +1.67%  "
$ws.Range("D10").Value = "'19.04"
$ws.Range("E10").Value = "'  +3.73%  "
$ws.Range("D12").Value = "'1.862.85"
$ws.Range("E12").Value = "'  +1.99%  "
$ws.Range("D13").Value = "'1.635.95"
$ws.Range("E13").Value = "'  +1.92%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "'  +1.44%  "
$ws.Range("E15").Value = "'  +2.49%  "
$ws.Range("D16").Value = "'26.671.52"
$ws.Range("E16").Value = "'  +1.87%  "
$ws.Range("D17").Value = "'63.00"
$ws.Range("E17").Value = "'  +1.86%  "
$ws.Range("E18").Value = "'  +1.84%  "
$ws.Range("E19").Value = "'  -0.01%  "
$ws.Range("D20").Value = "'208.61"
$ws.Range("E20").Value = "'  +3.99%  "
$ws.Range("E21").Value = "'  +0.82%  "
$ws.Range("D22").Value = "'9.39"
$ws.Range("E22").Value = "'  +1.22%  "
$ws.Range("E23").Value = "'  +2.91%  "
$ws.Range("D24").Value = "'1.90"
$ws.Range("E24").Value = "'  +2.21%  "
$ws.Range("D25").Value = "'146.44"
$ws.Range("E25").Value = "'  +1.48%  "
$ws.Range("E26").Value = "'  -0.04%  "
$ws.Range("E27").Value = "'  -0.81%  "
$ws.Range("D28").Value = "'6.73"
$ws.Range("E28").Value = "'  +2.72%  "
$ws.Range("E29").Value = "'  +1.22%  "
$ws.Range("E30").Value = "'  +5.74%  "
$ws.Range("E31").Value = "'  -0.44%  "
$ws.Range("E32").Value = "'  +1.19%  "
$ws.Range("E33").Value = "'  +0.99%  "
$ws.Range("E34").Value = "'  +1.50%  "
$ws.Range("E35").Value = "'  +0.82%  "
$ws.Range("D36").Value = "'1.169.07"
$ws.Range("E36").Value = "'  +0.55%  "
$ws.Range("E37").Value = "'  -0.67%  "
$ws.Range("E38").Value = "'  +2.55%  "
$ws.Range("E39").Value = "'  -0.01%  "
$ws.Range("E40").Value = "'  +0.36%  "
$ws.Range("E41").Value = "'  +1.56%  "
$ws.Range("E42").Value = "'  +1.84%  "
$ws.Range("D43").Value = "'5.37"
$ws.Range("E43").Value = "'  +1.34%  "
$ws.Range("D44").Value = "'1.773.23"
$ws.Range("E44").Value = "'  +2.01%  "
$ws.Range("D45").Value = "'92.30"
$ws.Range("E46").Value = "'  +1.84%  "
$ws.Range("E47").Value = "'  -2.21%  "
$ws.Range("D48").Value = "'54.69"
$ws.Range("E48").Value = "'  +1.16%  "
$ws.Range("E49").Value = "'  +1.62%  "
$ws.Range("E50").Value = "'  +0.81%  "
$ws.Range("D51").Value = "'7.52"
$ws.Range("E51").Value = "'  +4.47%  "
